$wb = $excel.ActiveWorkbook

# --- Rename the first sheet and add a second one -------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "OBJECT_DATA"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Planilha2"

# --- Fill in the decoded-object bytes (row 4) and their hex form (row 5) -
$ws2.Range("B4").Value = "dec"
$ws2.Range("C4").Value = 144
$ws2.Range("D4").Value = 57
$ws2.Range("E4").Value = 0
$ws2.Range("F4").Value = 31
$ws2.Range("G4").Value = 227
$ws2.Range("H4").Value = 252
$ws2.Range("I4").Value = 36
$ws2.Range("J4").Value = 166

$ws2.Range("B5").Value = "hex"
$ws2.Range("C5:J5").FormulaR1C1 = "=DEC2HEX(R[-1]C)"

# --- Match the page margins used on the rest of the workbook -------------
$ps2 = $ws2.PageSetup
$ps2.LeftMargin = 36.850393728
$ps2.RightMargin = 36.850393728
$ps2.TopMargin = 56.692913399999995
$ps2.BottomMargin = 56.692913399999995
$ps2.HeaderMargin = 22.67716464
$ps2.FooterMargin = 22.67716464

# --- Restore the view state on both sheets --------------------------------
$ws1.Activate()
$ws1.Range("J5:P6").Select()

$ws2.Activate()
$ws2.Range("I12").Select()
